$d = $word.ActiveDocument
$wmain = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Edit 1: "Como um usuário..." paragraph -> add jc=both and split the run so
# "Eu" is wrapped with proofErr gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Como um usuário do sistema em busca de autoaperfeiçoamento e progresso, Eu desejo uma visualização clara e abrangente do meu progresso ao longo da jornada do herói, Para ter uma compreensão detalhada das tarefas concluídas e em andamento, mantendo-me motivado e focado em meus objetivos.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $p = $find.Parent.Paragraphs(1)
    $xml = @"
<w:p xmlns:w="$wmain">
  <w:pPr>
    <w:spacing w:after="160" w:line="259" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Como um usuário do sistema em busca de autoaperfeiçoamento e progresso, </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Eu</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> desejo uma visualização clara e abrangente do meu progresso ao longo da jornada do herói, Para ter uma compreensão detalhada das tarefas concluídas e em andamento, mantendo-me motivado e focado em meus objetivos.</w:t>
  </w:r>
</w:p>
"@
    $p.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Edit 2: the empty paragraph right after "A interface deve incluir..." is
# replaced by three new list paragraphs (Interface Lúdica de Progresso: / two
# sub-bullets).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("A interface deve incluir elementos gráficos, como barras de conclusão ou gráficos de progresso, para representar visualmente o avanço do usuário.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $find.Parent.Paragraphs(1)
    $nextPara = $para.Next()
    $xml = @"
<w:p xmlns:w="$wmain">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="26"/>
    </w:numPr>
    <w:spacing w:after="160" w:line="259" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Interface </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>L</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">údica de </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>P</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>rogresso</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$wmain">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="26"/>
    </w:numPr>
    <w:spacing w:after="160" w:line="259" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>O sistema deve apresentar ao usuário</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> a porcentagem de conclusão de sua jornada</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="$wmain">
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="26"/>
    </w:numPr>
    <w:spacing w:after="160" w:line="259" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>O usuário deve visualizar em tempo real o percurso do seu personagem durante a conclusão das ta</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>refas.</w:t>
  </w:r>
</w:p>
"@
    $nextPara.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Edit 3: "Exibição do progresso geral na jornada através de uma porcentagem
# de conclusão." -> "... através de uma porcentagem a medida que as tarefas
# são concluídas." (with proofErr around "a").
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Exibição do progresso geral na jornada através de uma porcentagem de conclusão.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $p = $find.Parent.Paragraphs(1)
    $xml = @"
<w:p xmlns:w="$wmain">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="19"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Exibição do progresso geral na jornada</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> através de uma porcentagem </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>a</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> medida que as tarefas são concluídas.</w:t>
  </w:r>
</w:p>
"@
    $p.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Edit 4: "Pontuação total acumulada com a conclusão das atividades." ->
# "Porcentagem total acumulada com a conclusão das atividades." (word split
# into "Po" + "rcentagem" runs), and the following empty paragraph
# (ind left=708) is removed.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Pontuação total acumulada com a conclusão das atividades.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $p = $find.Parent.Paragraphs(1)
    $nextPara = $p.Next()
    $xml = @"
<w:p xmlns:w="$wmain">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="19"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Po</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>rcentagem</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">total </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>acumulada</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> com a conclusão das atividades</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
"@
    $p.Range.InsertXML($xml)
    $nextPara.Range.Delete()
}
